$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "62.591.77"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.34%  "
$cell.Style = "Normal"
# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.338.25"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.32%  "
$cell.Style = "Normal"
# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "560.03"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.36%  "
$cell.Style = "Normal"
# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "151.88"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.56%  "
$cell.Style = "Normal"
# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.14%  "
$cell.Style = "Normal"
# Row 8
$cell = $ws.Cells.Item(8, 2)
$cell.NumberFormat = "@"
$cell.Value = "LidoStakedEther"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.340.61"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.19%  "
$cell.Style = "Normal"
# Row 9
$cell = $ws.Cells.Item(9, 2)
$cell.NumberFormat = "@"
$cell.Value = "XRP"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.532"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.11%  "
$cell.Style = "Normal"
# Row 10
$cell = $ws.Cells.Item(10, 2)
$cell.NumberFormat = "@"
$cell.Value = "Toncoin"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.39"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.23%  "
$cell.Style = "Normal"
# Row 11
$cell = $ws.Cells.Item(11, 2)
$cell.NumberFormat = "@"
$cell.Value = "Dogecoin"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.117"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.69%  "
$cell.Style = "Normal"
# Row 12
$cell = $ws.Cells.Item(12, 2)
$cell.NumberFormat = "@"
$cell.Value = "Cardano"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.434"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.03%  "
$cell.Style = "Normal"
# Row 13
$cell = $ws.Cells.Item(13, 2)
$cell.NumberFormat = "@"
$cell.Value = "WrappedliquidstakedEther2.0"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.915.28"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.31%  "
$cell.Style = "Normal"
# Row 14
$cell = $ws.Cells.Item(14, 2)
$cell.NumberFormat = "@"
$cell.Value = "TRON"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.138"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.22%  "
$cell.Style = "Normal"
# Row 15
$cell = $ws.Cells.Item(15, 2)
$cell.NumberFormat = "@"
$cell.Value = "Avalanche"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.81"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.93%  "
$cell.Style = "Normal"
# Row 16
$cell = $ws.Cells.Item(16, 2)
$cell.NumberFormat = "@"
$cell.Value = "ShibaInu"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000178"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.06%  "
$cell.Style = "Normal"
# Row 17
$cell = $ws.Cells.Item(17, 2)
$cell.NumberFormat = "@"
$cell.Value = "WrappedBTC"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "62.609.26"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.22%  "
$cell.Style = "Normal"
# Row 18
$cell = $ws.Cells.Item(18, 2)
$cell.NumberFormat = "@"
$cell.Value = "WrappedEther"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.362.90"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.09%  "
$cell.Style = "Normal"
# Row 19
$cell = $ws.Cells.Item(19, 2)
$cell.NumberFormat = "@"
$cell.Value = "Polkadot"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.33"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.29%  "
$cell.Style = "Normal"
# Row 20
$cell = $ws.Cells.Item(20, 2)
$cell.NumberFormat = "@"
$cell.Value = "Chainlink"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "13.77"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.47%  "
$cell.Style = "Normal"
# Row 21
$cell = $ws.Cells.Item(21, 2)
$cell.NumberFormat = "@"
$cell.Value = "Uniswap"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.35"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.64%  "
$cell.Style = "Normal"
# Row 22
$cell = $ws.Cells.Item(22, 2)
$cell.NumberFormat = "@"
$cell.Value = "BitcoinCash"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "383.15"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.09%  "
$cell.Style = "Normal"
# Row 23
$cell = $ws.Cells.Item(23, 2)
$cell.NumberFormat = "@"
$cell.Value = "Dai"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.08%  "
$cell.Style = "Normal"
# Row 24
$cell = $ws.Cells.Item(24, 2)
$cell.NumberFormat = "@"
$cell.Value = "Polygon"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.533"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.59%  "
$cell.Style = "Normal"
# Row 25
$cell = $ws.Cells.Item(25, 2)
$cell.NumberFormat = "@"
$cell.Value = "Litecoin"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "70.00"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.16%  "
$cell.Style = "Normal"
# Row 26
$cell = $ws.Cells.Item(26, 2)
$cell.NumberFormat = "@"
$cell.Value = "Kaspa"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.178"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.10%  "
$cell.Style = "Normal"
# Row 27
$cell = $ws.Cells.Item(27, 2)
$cell.NumberFormat = "@"
$cell.Value = "InternetComputer(DFINITY)"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.91"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.62%  "
$cell.Style = "Normal"
# Row 28
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.08%  "
$cell.Style = "Normal"
# Row 29
$cell = $ws.Cells.Item(29, 2)
$cell.NumberFormat = "@"
$cell.Value = "PEPE"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0944"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.37%  "
$cell.Style = "Normal"
# Row 30
$cell = $ws.Cells.Item(30, 2)
$cell.NumberFormat = "@"
$cell.Value = "RenderToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.53"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.49%  "
$cell.Style = "Normal"
# Row 31
$cell = $ws.Cells.Item(31, 2)
$cell.NumberFormat = "@"
$cell.Value = "PancakeSwap"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.98"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.86%  "
$cell.Style = "Normal"
# Row 32
$cell = $ws.Cells.Item(32, 2)
$cell.NumberFormat = "@"
$cell.Value = "NEARProtocol"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.56"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.55%  "
$cell.Style = "Normal"
# Row 33
$cell = $ws.Cells.Item(33, 2)
$cell.NumberFormat = "@"
$cell.Value = "EthereumClassic"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "22.84"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.27%  "
$cell.Style = "Normal"
# Row 34
$cell = $ws.Cells.Item(34, 2)
$cell.NumberFormat = "@"
$cell.Value = "Fetch.AI"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.30"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +7.51%  "
$cell.Style = "Normal"
# Row 35
$cell = $ws.Cells.Item(35, 2)
$cell.NumberFormat = "@"
$cell.Value = "Aptos"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.68"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.Style = "Normal"
# Row 36
$cell = $ws.Cells.Item(36, 2)
$cell.NumberFormat = "@"
$cell.Value = "Monero"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "159.95"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.83%  "
$cell.Style = "Normal"
# Row 37
$cell = $ws.Cells.Item(37, 2)
$cell.NumberFormat = "@"
$cell.Value = "ImmutableX"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.47"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +8.19%  "
$cell.Style = "Normal"
# Row 38
$cell = $ws.Cells.Item(38, 2)
$cell.NumberFormat = "@"
$cell.Value = "Stacks"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.89"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +12.98%  "
$cell.Style = "Normal"
# Row 39
$cell = $ws.Cells.Item(39, 2)
$cell.NumberFormat = "@"
$cell.Value = "EnergySwap"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.81"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.07%  "
$cell.Style = "Normal"
# Row 40
$cell = $ws.Cells.Item(40, 2)
$cell.NumberFormat = "@"
$cell.Value = "Hedera"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0737"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.65%  "
$cell.Style = "Normal"
# Row 41
$cell = $ws.Cells.Item(41, 2)
$cell.NumberFormat = "@"
$cell.Value = "Maker"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.793.68"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.16%  "
$cell.Style = "Normal"
# Row 42
$cell = $ws.Cells.Item(42, 2)
$cell.NumberFormat = "@"
$cell.Value = "VeChain"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0312"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +6.19%  "
$cell.Style = "Normal"
# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "40.48"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.62%  "
$cell.Style = "Normal"
# Row 44
$cell = $ws.Cells.Item(44, 2)
$cell.NumberFormat = "@"
$cell.Value = "Mantle"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.742"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.38%  "
$cell.Style = "Normal"
# Row 45
$cell = $ws.Cells.Item(45, 2)
$cell.NumberFormat = "@"
$cell.Value = "Filecoin"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.24"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.09%  "
$cell.Style = "Normal"
# Row 46
$cell = $ws.Cells.Item(46, 2)
$cell.NumberFormat = "@"
$cell.Value = "ONDO"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.03"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.50%  "
$cell.Style = "Normal"
# Row 47
$cell = $ws.Cells.Item(47, 2)
$cell.NumberFormat = "@"
$cell.Value = "RenzoRestakedETH"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.385.21"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.33%  "
$cell.Style = "Normal"
# Row 48
$cell = $ws.Cells.Item(48, 2)
$cell.NumberFormat = "@"
$cell.Value = "InjectiveProtocol"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "21.84"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +6.01%  "
$cell.Style = "Normal"
# Row 49
$cell = $ws.Cells.Item(49, 2)
$cell.NumberFormat = "@"
$cell.Value = "Stellar"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.00%  "
$cell.Style = "Normal"
# Row 50
$cell = $ws.Cells.Item(50, 2)
$cell.NumberFormat = "@"
$cell.Value = "Cosmos"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.29"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.03%  "
$cell.Style = "Normal"
# Row 51
$cell = $ws.Cells.Item(51, 2)
$cell.NumberFormat = "@"
$cell.Value = "Bittensor"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "286.41"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +6.33%  "
$cell.Style = "Normal"
